$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 11
$ws.Range("A11").Value = "AD07 beregn indtjeningsbidrag"
$ws.Range("C11").Value = 43888
$ws.Range("D11").Value = 0.36458333333333331
$ws.Range("E11").Value = 0.40138888888888885
$ws.Range("F11").Value = "1 time"

# Row 12
$ws.Range("A12").Value = "Review UC02"
$ws.Range("B12").Value = "reviewer"
$ws.Range("C12").Value = 43888
$ws.Range("D12").Value = 0.40972222222222227
$ws.Range("E12").Value = 0.41666666666666669
$ws.Range("F12").Value = "10 til 20 min"

# Row 13
$ws.Range("A13").Value = "review DOM02"
$ws.Range("B13").Value = "reviewer"
$ws.Range("C13").Value = 43888
$ws.Range("D13").Value = 0.41666666666666669
$ws.Range("E13").Value = 0.42499999999999999
$ws.Range("F13").Value = "10 til 20 min"

# Row 14
$ws.Range("A14").Value = "Rette OC06"
$ws.Range("B14").Value = "System Analyst "
$ws.Range("C14").Value = 43889
$ws.Range("D14").Value = 0.36458333333333331
$ws.Range("E14").Value = 0.37152777777777773
$ws.Range("F14").Value = "10min"

# Row 16 (shared string "ATD06a" was registered before "ATD07b" in the
# original authoring session, so write this row's label first).
$ws.Range("A16").Value = "ATD06a"
$ws.Range("C16").Value = 43889
$ws.Range("D16").Value = 0.40625
$ws.Range("E16").Value = 0.42708333333333331
$ws.Range("F16").Value = "30 min"

# Row 15
$ws.Range("A15").Value = "ATD07b"
$ws.Range("C15").Value = 43889
$ws.Range("D15").Value = 0.375
$ws.Range("E15").Value = 0.39583333333333331
$ws.Range("F15").Value = "30 min"

# Row 17
$ws.Range("A17").Value = "DD07"
$ws.Range("C17").Value = 43889
$ws.Range("D17").Value = 0.43055555555555558
$ws.Range("E17").Value = 0.44791666666666669
$ws.Range("F17").Value = "20 min"

# Update the active selection to C17, matching the authored workbook state.
[void]$ws.Range("C17").Select()
